$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that were repulled/recalculated.
$ws.Range("F3").Value = 4
$ws.Range("F5").Value = -1
$ws.Range("F9").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = -2
$ws.Range("F31").Value = 3
$ws.Range("F39").Value = -1
$ws.Range("F45").Value = 3
$ws.Range("F50").Value = -1
$ws.Range("F52").Value = 4
$ws.Range("F55").Value = 2
$ws.Range("F56").Value = 2
$ws.Range("F62").Value = 6
